# Rebuild Sheet1 data: the export now has THREE data columns (Black_US,
# Black_Africa, White_US) instead of two, plus a Comments column shifted to D.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out the old 2-column layout (A/B = weights, C = comments) before
# laying down the corrected 3-column layout (A/B/C = weights, D = comments).
$ws.UsedRange.Clear()

# Row 1
$ws.Range('A1').Value = 'Black_US'
$ws.Range('B1').Value = 'Black_Africa'
$ws.Range('C1').Value = 'White_US'
$ws.Range('D1').Value = 'Comments'

# Row 2
$ws.Range('A2').Value = 2566
$ws.Range('B2').Value = 2790
$ws.Range('C2').Value = 4429
$ws.Range('D2').Value = 'Data representative of results published in:'

# Row 3
$ws.Range('A3').Value = 3068
$ws.Range('B3').Value = 3283
$ws.Range('C3').Value = 3191
$ws.Range('D3').Value = 'David, Richard J., & Collins, James W.'

# Row 4
$ws.Range('A4').Value = 2931
$ws.Range('B4').Value = 3148
$ws.Range('C4').Value = 3712
$ws.Range('D4').Value = 'Differing birth weight among infants of U.S.-born'

# Row 5
$ws.Range('A5').Value = 2882
$ws.Range('B5').Value = 3101
$ws.Range('C5').Value = 3399
$ws.Range('D5').Value = 'blacks, African-born blacks, and U.S.-born whites.'

# Row 6
$ws.Range('A6').Value = 2022
$ws.Range('B6').Value = 2257
$ws.Range('C6').Value = 2638
$ws.Range('D6').Value = 'New England Journal of Medicine.'

# Row 7
$ws.Range('A7').Value = 2073
$ws.Range('B7').Value = 2307
$ws.Range('C7').Value = 3946
$ws.Range('D7').Value = '337(17):1209-1214.'

# Row 8
$ws.Range('A8').Value = 2907
$ws.Range('B8').Value = 3125
$ws.Range('C8').Value = 3173
$ws.Range('D8').Value = ''''

# Row 9
$ws.Range('A9').Value = 4028
$ws.Range('B9').Value = 4225
$ws.Range('C9').Value = 2926
$ws.Range('D9').Value = 'The data are representive of birthweights (in grams)'

# Row 10
$ws.Range('A10').Value = 2982
$ws.Range('B10').Value = 3199
$ws.Range('C10').Value = 2303
$ws.Range('D10').Value = 'of children born in Illinois to mothers who fall into'

# Row 11
$ws.Range('A11').Value = 2907
$ws.Range('B11').Value = 3125
$ws.Range('C11').Value = 3885
$ws.Range('D11').Value = 'one of the following categories:'

# Row 12
$ws.Range('A12').Value = 2893
$ws.Range('B12').Value = 3112
$ws.Range('C12').Value = 3208
$ws.Range('D12').Value = '(1) Black, born in the United States (Black_US)'

# Row 13
$ws.Range('A13').Value = 2422
$ws.Range('B13').Value = 2649
$ws.Range('C13').Value = 2969
$ws.Range('D13').Value = '(2) Black, born in Africa (Black_Africa), or'

# Row 14
$ws.Range('A14').Value = 3910
$ws.Range('B14').Value = 4109
$ws.Range('C14').Value = 2948
$ws.Range('D14').Value = '(3) White, born in the United States (White_US).'

# Row 15
$ws.Range('A15').Value = 2588
$ws.Range('B15').Value = 2812
$ws.Range('C15').Value = 2270
$ws.Range('D15').Value = ''''

# Row 16
$ws.Range('A16').Value = 2832
$ws.Range('B16').Value = 3051
$ws.Range('C16').Value = 3172
$ws.Range('D16').Value = ''''

# Row 17
$ws.Range('A17').Value = 2063
$ws.Range('B17').Value = 2297
$ws.Range('C17').Value = 2318
$ws.Range('D17').Value = ''''

# Row 18
$ws.Range('A18').Value = 2213
$ws.Range('B18').Value = 2444
$ws.Range('C18').Value = 2456
$ws.Range('D18').Value = ''''

# Row 19
$ws.Range('A19').Value = 3672
$ws.Range('B19').Value = 3875
$ws.Range('C19').Value = 3661
$ws.Range('D19').Value = ''''

# Row 20
$ws.Range('A20').Value = 3512
$ws.Range('B20').Value = 3718
$ws.Range('C20').Value = 3854
$ws.Range('D20').Value = ''''

# Row 21
$ws.Range('A21').Value = 3425
$ws.Range('B21').Value = 3633
$ws.Range('C21').Value = 3122
$ws.Range('D21').Value = ''''

# Row 22
$ws.Range('A22').Value = 3097
$ws.Range('B22').Value = 3312
$ws.Range('C22').Value = 3666
$ws.Range('D22').Value = ''''

# Row 23
$ws.Range('A23').Value = 2583
$ws.Range('B23').Value = 2807
$ws.Range('C23').Value = 4414
$ws.Range('D23').Value = ''''

# Row 24
$ws.Range('A24').Value = 3961
$ws.Range('B24').Value = 4159
$ws.Range('C24').Value = 3490
$ws.Range('D24').Value = ''''

# Row 25
$ws.Range('A25').Value = 2936
$ws.Range('B25').Value = 3154
$ws.Range('C25').Value = 3871
$ws.Range('D25').Value = ''''

# Row 26
$ws.Range('A26').Value = 3081
$ws.Range('B26').Value = 3296
$ws.Range('C26').Value = 2679
$ws.Range('D26').Value = ''''

# Row 27
$ws.Range('A27').Value = 3873
$ws.Range('B27').Value = 4073
$ws.Range('C27').Value = 2850
$ws.Range('D27').Value = ''''

# Row 28
$ws.Range('A28').Value = 2966
$ws.Range('B28').Value = 3183
$ws.Range('C28').Value = 2852
$ws.Range('D28').Value = ''''

# Row 29
$ws.Range('A29').Value = 2906
$ws.Range('B29').Value = 3124
$ws.Range('C29').Value = 3316
$ws.Range('D29').Value = ''''

# Row 30
$ws.Range('A30').Value = 3489
$ws.Range('B30').Value = 3696
$ws.Range('C30').Value = 3596
$ws.Range('D30').Value = ''''

# Row 31
$ws.Range('A31').Value = 2525
$ws.Range('B31').Value = 2751
$ws.Range('C31').Value = 2719
$ws.Range('D31').Value = ''''

# Row 32
$ws.Range('A32').Value = 3087
$ws.Range('B32').Value = 3302
$ws.Range('C32').Value = 4448
$ws.Range('D32').Value = ''''

# Row 33
$ws.Range('A33').Value = 3697
$ws.Range('B33').Value = 3900
$ws.Range('C33').Value = 3043
$ws.Range('D33').Value = ''''

# Row 34
$ws.Range('A34').Value = 2293
$ws.Range('B34').Value = 2523
$ws.Range('C34').Value = 2709
$ws.Range('D34').Value = ''''

# Row 35
$ws.Range('A35').Value = 2928
$ws.Range('B35').Value = 3146
$ws.Range('C35').Value = 3695
$ws.Range('D35').Value = ''''

# Row 36
$ws.Range('A36').Value = 2641
$ws.Range('B36').Value = 2864
$ws.Range('C36').Value = 3583
$ws.Range('D36').Value = ''''

# Row 37
$ws.Range('B37').Value = 3126
$ws.Range('C37').Value = 2867
$ws.Range('D37').Value = ''''

# Row 38
$ws.Range('B38').Value = 2661
$ws.Range('C38').Value = 4056
$ws.Range('D38').Value = ''''

# Row 39
$ws.Range('B39').Value = 3951
$ws.Range('C39').Value = 3342
$ws.Range('D39').Value = ''''

# Row 40
$ws.Range('B40').Value = 4068
$ws.Range('C40').Value = 3124
$ws.Range('D40').Value = ''''

# Row 41
$ws.Range('C41').Value = 4281
$ws.Range('D41').Value = ''''

# Row 42
$ws.Range('C42').Value = 3839
$ws.Range('D42').Value = ''''

# Row 43
$ws.Range('C43').Value = 3458
$ws.Range('D43').Value = ''''

# Row 44
$ws.Range('C44').Value = 3931
$ws.Range('D44').Value = ''''

# Row 45
$ws.Range('C45').Value = 4322
$ws.Range('D45').Value = ''''

# Restore the active cell to the default (A1) now that the old C1 selection
# (the stray comments header) no longer reflects the sheet layout.
$ws.Range("A1").Select() | Out-Null

